$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 312: Calameño/Extra (old, 2021-01-13) -> Tuna/Extra (new week, 2023-01-13) ---
$ws.Range("D312").Value2 = 44939
$ws.Range("H312").Value2 = "Tuna"
$ws.Range("J312").Value2 = 500
$ws.Range("K312").Value2 = 1400
$ws.Range("L312").Value2 = 1400
$ws.Range("M312").Value2 = 1400
$ws.Range("P312").Value2 = 1400

# --- Row 313: Calameño/Primera (old) -> Tuna/Primera (new week) ---
$ws.Range("D313").Value2 = 44939
$ws.Range("H313").Value2 = "Tuna"
$ws.Range("J313").Value2 = 1000
$ws.Range("K313").Value2 = 1100
$ws.Range("L313").Value2 = 1100
$ws.Range("M313").Value2 = 1100
$ws.Range("P313").Value2 = 1100

# --- Row 314: Calameño/Segunda (old) -> Tuna/Segunda (new week) ---
$ws.Range("D314").Value2 = 44939
$ws.Range("H314").Value2 = "Tuna"
$ws.Range("J314").Value2 = 500
$ws.Range("K314").Value2 = 800
$ws.Range("L314").Value2 = 800
$ws.Range("M314").Value2 = 800
$ws.Range("P314").Value2 = 800

# --- Rows 315-317 previously held "Tuna" (2021-01-13); they now hold the
# --- "Calameño" data that used to live in rows 312-314 (only H changes). ---
$ws.Range("H315").Value2 = "Calameño"
$ws.Range("H316").Value2 = "Calameño"
$ws.Range("H317").Value2 = "Calameño"

# --- New rows 318-320: the "Tuna" 2021-01-13 rows that used to be 315-317,
# --- now appended at the end of the table. ---
$ws.Range("A318").Value2 = 11
$ws.Range("B318").Value2 = "Vega Monumental Concepción"
$ws.Range("C318").Value2 = "Bíobío"
$ws.Range("D318").NumberFormat = $ws.Range("D317").NumberFormat
$ws.Range("D318").Value2 = 44209
$ws.Range("E318").Value2 = 8
$ws.Range("F318").Value2 = 100112027
$ws.Range("G318").Value2 = "Melón"
$ws.Range("H318").Value2 = "Tuna"
$ws.Range("I318").Value2 = "Extra"
$ws.Range("J318").Value2 = 400
$ws.Range("K318").Value2 = 1000
$ws.Range("L318").Value2 = 1000
$ws.Range("M318").Value2 = 1000
$ws.Range("N318").Value2 = "`$/unidad"
$ws.Range("O318").Value2 = "Región de O'Higgins"
$ws.Range("P318").Value2 = 1000
$ws.Range("Q318").Value2 = 1
$ws.Range("R318").Value2 = "Hortaliza"

$ws.Range("A319").Value2 = 11
$ws.Range("B319").Value2 = "Vega Monumental Concepción"
$ws.Range("C319").Value2 = "Bíobío"
$ws.Range("D319").NumberFormat = $ws.Range("D317").NumberFormat
$ws.Range("D319").Value2 = 44209
$ws.Range("E319").Value2 = 8
$ws.Range("F319").Value2 = 100112027
$ws.Range("G319").Value2 = "Melón"
$ws.Range("H319").Value2 = "Tuna"
$ws.Range("I319").Value2 = "Primera"
$ws.Range("J319").Value2 = 400
$ws.Range("K319").Value2 = 800
$ws.Range("L319").Value2 = 800
$ws.Range("M319").Value2 = 800
$ws.Range("N319").Value2 = "`$/unidad"
$ws.Range("O319").Value2 = "Región de O'Higgins"
$ws.Range("P319").Value2 = 800
$ws.Range("Q319").Value2 = 1
$ws.Range("R319").Value2 = "Hortaliza"

$ws.Range("A320").Value2 = 11
$ws.Range("B320").Value2 = "Vega Monumental Concepción"
$ws.Range("C320").Value2 = "Bíobío"
$ws.Range("D320").NumberFormat = $ws.Range("D317").NumberFormat
$ws.Range("D320").Value2 = 44209
$ws.Range("E320").Value2 = 8
$ws.Range("F320").Value2 = 100112027
$ws.Range("G320").Value2 = "Melón"
$ws.Range("H320").Value2 = "Tuna"
$ws.Range("I320").Value2 = "Segunda"
$ws.Range("J320").Value2 = 400
$ws.Range("K320").Value2 = 600
$ws.Range("L320").Value2 = 600
$ws.Range("M320").Value2 = 600
$ws.Range("N320").Value2 = "`$/unidad"
$ws.Range("O320").Value2 = "Región de O'Higgins"
$ws.Range("P320").Value2 = 600
$ws.Range("Q320").Value2 = 1
$ws.Range("R320").Value2 = "Hortaliza"
